# Auto-generated Excel COM-interop script applying the scheduled-runner data refresh
# to the Famfrit_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 504.08334
$ws.Range("J13").Value = 549.5
$ws.Range("L13").Value = 549.5
$ws.Range("N13").Value = -887.5
$ws.Range("H76").Value = 10754.5
$ws.Range("I76").Value = 18084.428
$ws.Range("J76").Value = 6090
$ws.Range("K76").Value = 18084.428
$ws.Range("L76").Value = 6090
$ws.Range("M76").Value = -17769.428
$ws.Range("N76").Value = -6720
$ws.Range("H79").Value = 10754.5
$ws.Range("I79").Value = 18084.428
$ws.Range("J79").Value = 6090
$ws.Range("K79").Value = 18084.428
$ws.Range("L79").Value = 6090
$ws.Range("M79").Value = -16992.428
$ws.Range("N79").Value = -8274
$ws.Range("H98").Value = 1390.421
$ws.Range("I98").Value = 1410.6666
$ws.Range("J98").Value = 1026
$ws.Range("K98").Value = 1410.6666
$ws.Range("L98").Value = 1026
$ws.Range("M98").Value = 87.33339999999998
$ws.Range("N98").Value = -4022
$ws.Range("H100").Value = 3000.2
$ws.Range("I100").Value = 1517.3334
$ws.Range("J100").Value = 5224.5
$ws.Range("K100").Value = 1517.3334
$ws.Range("L100").Value = 5224.5
$ws.Range("M100").Value = -976.3334
$ws.Range("N100").Value = -6306.5
$ws.Range("H112").Value = 5683817
$ws.Range("I112").Value = 1000
$ws.Range("K112").Value = 3000
$ws.Range("M112").Value = -1892
$ws.Range("H122").Value = 1390.421
$ws.Range("I122").Value = 1410.6666
$ws.Range("J122").Value = 1026
$ws.Range("K122").Value = 4231.9998
$ws.Range("L122").Value = 3078
$ws.Range("M122").Value = -1781.9998
$ws.Range("N122").Value = -7978
$ws.Range("H132").Value = 2676.4856
$ws.Range("I132").Value = 2043
$ws.Range("K132").Value = 6129
$ws.Range("M132").Value = -3599
$ws.Range("H137").Value = 1976.5238
$ws.Range("I137").Value = 2080.3333
$ws.Range("J137").Value = 1838.1111
$ws.Range("K137").Value = 6240.999899999999
$ws.Range("L137").Value = 5514.3333
$ws.Range("M137").Value = -3690.999899999999
$ws.Range("N137").Value = -10614.3333
$ws.Range("H138").Value = 15158003
$ws.Range("J138").Value = 22230908
$ws.Range("L138").Value = 66692724
$ws.Range("N138").Value = -66703004
$ws.Range("H141").Value = 1817.4736
$ws.Range("I141").Value = 1817.4736
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5452.4208
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -272.4207999999999
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1795.5238
$ws.Range("I45").Value = 1554.5883
$ws.Range("K45").Value = 1554.5883
$ws.Range("M45").Value = -1177.5883
$ws.Range("H74").Value = 23282446
$ws.Range("I74").Value = 23282446
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 23282446
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -23281572
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 23282446
$ws.Range("I77").Value = 23282446
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 116412230
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -116407862
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 26355380
$ws.Range("I132").Value = 2752.4688
$ws.Range("J132").Value = 166902740
$ws.Range("K132").Value = 8257.4064
$ws.Range("L132").Value = 500708220
$ws.Range("M132").Value = -5727.4064
$ws.Range("N132").Value = -500713280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8527.866
$ws.Range("I105").Value = 9784.833000000001
$ws.Range("K105").Value = 9784.833000000001
$ws.Range("M105").Value = -8037.833000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3856.8708
$ws.Range("I31").Value = 2251.6667
$ws.Range("J31").Value = 8105.9414
$ws.Range("K31").Value = 2251.6667
$ws.Range("L31").Value = 8105.9414
$ws.Range("M31").Value = -1956.6667
$ws.Range("N31").Value = -8695.9414
$ws.Range("H34").Value = 3856.8708
$ws.Range("I34").Value = 2251.6667
$ws.Range("J34").Value = 8105.9414
$ws.Range("K34").Value = 2251.6667
$ws.Range("L34").Value = 8105.9414
$ws.Range("M34").Value = -2049.6667
$ws.Range("N34").Value = -8509.9414
$ws.Range("H107").Value = 1935.9445
$ws.Range("I107").Value = 788.6923
$ws.Range("K107").Value = 788.6923
$ws.Range("M107").Value = 1131.3077

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3461
$ws.Range("J113").Value = 3856.4546
$ws.Range("L113").Value = 11569.3638
$ws.Range("N113").Value = -15909.3638
$ws.Range("H122").Value = 1437
$ws.Range("I122").Value = 799.6667
$ws.Range("K122").Value = 7197.0003
$ws.Range("M122").Value = -4747.0003
$ws.Range("H132").Value = 4227
$ws.Range("I132").Value = 2302.6667
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 20724.0003
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -18194.0003
$ws.Range("N132").Value = -95060

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2606.6956
$ws.Range("I102").Value = 547.3333
$ws.Range("K102").Value = 547.3333
$ws.Range("M102").Value = 1074.6667
$ws.Range("H132").Value = 3052.6843
$ws.Range("I132").Value = 3143.1428
$ws.Range("J132").Value = 2799.4
$ws.Range("K132").Value = 9429.428400000001
$ws.Range("L132").Value = 8398.200000000001
$ws.Range("M132").Value = -6899.428400000001
$ws.Range("N132").Value = -13458.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2118.7742
$ws.Range("I40").Value = 2052.72
$ws.Range("K40").Value = 2052.72
$ws.Range("M40").Value = -1916.72
$ws.Range("H100").Value = 3376.7778
$ws.Range("I100").Value = 3038.3333
$ws.Range("K100").Value = 3038.3333
$ws.Range("M100").Value = -2497.3333
$ws.Range("H122").Value = 3665.1064
$ws.Range("I122").Value = 2509.84
$ws.Range("K122").Value = 7529.52
$ws.Range("M122").Value = -5079.52
$ws.Range("H136").Value = 2031.0312
$ws.Range("I136").Value = 1708.0834
$ws.Range("J136").Value = 2999.875
$ws.Range("K136").Value = 5124.2502
$ws.Range("L136").Value = 8999.625
$ws.Range("M136").Value = -2574.2502
$ws.Range("N136").Value = -14099.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6354.154
$ws.Range("I126").Value = 6145.8184
$ws.Range("K126").Value = 18437.4552
$ws.Range("M126").Value = -15967.4552
$ws.Range("H136").Value = 1120.25
$ws.Range("I136").Value = 807.8889
$ws.Range("J136").Value = 2525.875
$ws.Range("K136").Value = 2423.6667
$ws.Range("L136").Value = 7577.625
$ws.Range("M136").Value = 126.3332999999998
$ws.Range("N136").Value = -12677.625

